# TeamVOXSReport.xlsx - Sprint1 sheet updates:
#  - US02 (row 3) and US07 (row 8) completed: Status -> Done, Act Size/Act
#    Time filled in, Completed column marked "Yes".
#  - US06 (row 7) and US08 (row 9) reassigned from SR/VB to a new owner "XP".
#  - The stray "Completed" date on US01 (row 2, I2) is cleared.
#  - US09 (row 10) is wiped back to blank (only the existing cell shading remains).
#  - The active sheet/tab moves from "Burndown" to "Sprint1", with the
#    selection resting on H9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Row 2: the "Completed" date on US01 is removed entirely (cell + format).
$ws.Range("I2").Clear()

# Owner reassignments to "XP" (write these before any "Done"/"Yes" text below
# so new shared strings get created in the same order as the target file:
# XP, then Done, then Yes).
$ws.Range("C7").Value = "XP"   # US06: SR -> XP
$ws.Range("C9").Value = "XP"   # US08: VB -> XP

# Row 3 (US02): mark done, record actuals, flag as completed.
$ws.Range("D3").Value = "Done"
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = "Yes"

# Row 8 (US07): mark done, record actuals, flag as completed.
$ws.Range("D8").Value = "Done"
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 15
$ws.Range("I8").Value = "Yes"

# Row 10 (US09): clear the whole row - A:C keep their grey style but lose
# their content, D:F are cleared completely (content + format).
$ws.Range("A10:C10").ClearContents()
$ws.Range("D10:F10").Clear()

# View state: Sprint1 becomes the selected/active tab, cursor on H9.
$ws.Activate()
$ws.Range("H9").Select() | Out-Null
